# AfDD_2021_Annex_Table_Tab15.xlsx edit
#
# 1) The "Male life expectancy" / "Female life expectancy" header texts in
#    D2 / E2 were swapped (and, correspondingly, every data row's D/E
#    values were swapped so the numbers stay lined up under the correct
#    sex again).
# 2) A handful of blank "A" column label cells used in the sub-total /
#    regional-aggregate rows held a single space character; they are
#    cleared to a true empty string.
# 3) Cosmetic: the saved workbook window height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab15")

# --- 1) Swap Male / Female life-expectancy header text (row 2, cols D & E) ---
$headerD = $ws.Cells.Item(2, 4).Value2
$headerE = $ws.Cells.Item(2, 5).Value2
$ws.Cells.Item(2, 4).Value2 = $headerE
$ws.Cells.Item(2, 5).Value2 = $headerD

# --- Swap the D/E data values for every country / aggregate row (3-97) ---
for ($r = 3; $r -le 97; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
}

# --- 2) Clear the stray single-space label in column A of the subtotal rows ---
#    (scan the used range dynamically instead of a hard-coded row list)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq " ") {
        $cell.Value2 = ""
    }
}

# --- 3) Cosmetic saved-window size (best effort; may be a no-op headless) ---
$win = $wb.Windows.Item(1)
$win.Height = 11990
